$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily record row (2025/09/30) below the existing data.
$row = 37
$ws.Cells.Item($row, 1).Value = "'2025/09/30"
$ws.Cells.Item($row, 2).Value = "火"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 159
